$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '26.226.25'
$ws.Range('E2').Value = '  -2.25%  '
Set-TextValue $ws.Range('D3') '1.671.52'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  +0.33%  '
Set-TextValue $ws.Range('D5') '217.55'
$ws.Range('E5').Value = '  -1.45%  '
Set-TextValue $ws.Range('D6') '0.5119'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E7').Value = '  +0.25%  '
Set-TextValue $ws.Range('D8') '0.2655'
$ws.Range('E8').Value = '  +3.06%  '
Set-TextValue $ws.Range('D9') '0.06376'
$ws.Range('E9').Value = '  +2.88%  '
Set-TextValue $ws.Range('D10') '21.50'
$ws.Range('E10').Value = '  -2.49%  '
Set-TextValue $ws.Range('D11') '0.07395'
$ws.Range('E11').Value = '  +0.79%  '
Set-TextValue $ws.Range('D12') '1.675.83'
$ws.Range('E12').Value = '  -1.54%  '
Set-TextValue $ws.Range('D13') '4.543'
$ws.Range('E13').Value = '  +1.45%  '
Set-TextValue $ws.Range('D14') '0.5830'
$ws.Range('E14').Value = '  +0.42%  '
Set-TextValue $ws.Range('D15') '1.901.50'
$ws.Range('E15').Value = '  -1.83%  '
Set-TextValue $ws.Range('D16') '0.000008639'
$ws.Range('E16').Value = '  +5.39%  '
Set-TextValue $ws.Range('D17') '64.50'
$ws.Range('E17').Value = '  -1.55%  '
Set-TextValue $ws.Range('D18') '26.313.24'
$ws.Range('E18').Value = '  -2.16%  '
Set-TextValue $ws.Range('D19') '4.956'
$ws.Range('E19').Value = '  -0.97%  '
Set-TextValue $ws.Range('D21') '10.88'
$ws.Range('E21').Value = '  +2.12%  '
Set-TextValue $ws.Range('D22') '189.22'
$ws.Range('E22').Value = '  +1.56%  '
Set-TextValue $ws.Range('D23') '6.211'
$ws.Range('E23').Value = '  -0.88%  '
$ws.Range('E24').Value = '  +0.34%  '
Set-TextValue $ws.Range('D25') '144.12'
$ws.Range('E25').Value = '  +0.63%  '
Set-TextValue $ws.Range('D26') '7.643'
$ws.Range('E26').Value = '  +1.11%  '
Set-TextValue $ws.Range('D27') '0.1179'
$ws.Range('E27').Value = '  +2.13%  '
Set-TextValue $ws.Range('D28') '15.65'
$ws.Range('E28').Value = '  +2.90%  '
Set-TextValue $ws.Range('D29') '0.05988'
$ws.Range('E29').Value = '  +2.19%  '
Set-TextValue $ws.Range('D30') '1.284'
$ws.Range('E30').Value = '  -4.15%  '
Set-TextValue $ws.Range('D31') '1.326'
$ws.Range('E31').Value = '  -1.49%  '
Set-TextValue $ws.Range('D32') '3.526'
$ws.Range('E32').Value = '  +1.77%  '
Set-TextValue $ws.Range('D33') '3.522'
$ws.Range('E33').Value = '  +2.38%  '
$ws.Range('E34').Value = '  -0.07%  '
Set-TextValue $ws.Range('D35') '1.015'
$ws.Range('E35').Value = '  +2.78%  '
Set-TextValue $ws.Range('D36') '0.6026'
$ws.Range('E36').Value = '  +0.16%  '
Set-TextValue $ws.Range('D37') '2.376'
$ws.Range('E37').Value = '  -1.49%  '
Set-TextValue $ws.Range('D38') '2.659'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D39') '6.090'
$ws.Range('E39').Value = '  +3.60%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D40') '0.01613'
$ws.Range('E40').Value = '  +0.89%  '
Set-TextValue $ws.Range('D41') '1.080.82'
$ws.Range('E41').Value = '  -1.64%  '
Set-TextValue $ws.Range('D42') '0.8706'
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('E43').Value = '  +0.73%  '
Set-TextValue $ws.Range('D44') '100.29'
$ws.Range('E44').Value = '  +2.68%  '
Set-TextValue $ws.Range('D45') '1.822.89'
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('E46').Value = '  +7.32%  '
Set-TextValue $ws.Range('D47') '56.29'
$ws.Range('E47').Value = '  +0.29%  '
$ws.Range('E48').Value = '  -0.24%  '
Set-TextValue $ws.Range('D49') '8.080'
$ws.Range('E49').Value = '  +1.99%  '
Set-TextValue $ws.Range('D50') '0.05215'
$ws.Range('E50').Value = '  -0.64%  '
Set-TextValue $ws.Range('D51') '0.4300'
$ws.Range('E51').Value = '  -1.45%  '
